# Apply the update: append two new funding-agency rows (17 and 18) to the
# table on Sheet1, extending the table/dimension from A1:F93 to A1:F95.
#
# NOTE: the rows are written in the same order the original author typed
# them (Welch Foundation first, then the Rachadapisek Sompot Fund) even
# though the Rachadapisek row ends up physically above the Welch row, so
# that the resulting shared-string table ordering matches the source file.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 95: Welch Foundation (entered first)
$ws.Range("A95").Value = 18
$ws.Range("B95").Value = "Welch Foundation"
$ws.Range("D95").Value = "http://dx.doi.org/10.13039/100000928"
$ws.Range("E95").Value = "Y"

# Row 94: Rachadapisek Sompot Fund for Postdoctoral Fellowship, Chulalongkorn University (Thailand)
$ws.Range("A94").Value = 17
$ws.Range("B94").Value = "Rachadapisek Sompot Fund for Postdoctoral Fellowship, Chulalongkorn University (Thailand)"
$ws.Range("D94").Value = "http://dx.doi.org/10.13039/501100002873"
$ws.Range("E94").Value = "Y"

# Grow the worksheet table (ListObject) so the new rows are part of Table1
$table = $ws.ListObjects.Item("Table1")
$table.Resize($ws.Range("A1:F95"))

# Update the view so the currently visible selection matches the author's
# final cursor position after adding the rows.
$ws.Application.ActiveWindow.ScrollRow = 58
$ws.Range("D94").Select()
